$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Batch/Class" value in column H, row 1
$ws.Range("H1").Value = "21CSE"

# Update the active selection to the newly added cell, matching the
# workbook's saved cursor position after the edit.
$ws.Range("H1").Select()
